$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 46 with Sina mail entry
$ws.Range("A46").Value = "新浪邮箱"
$ws.Range("B46").Value = "bingona@sina.com"
$ws.Range("C46").Value = "Bingona1314"

# Add the actual hyperlink (mailto) on B46
$ws.Hyperlinks.Add($ws.Range("B46"), "mailto:bingona@sina.com", [Type]::Missing, [Type]::Missing, "bingona@sina.com")

# Re-apply the same cell style object as the other hyperlink cells (reuse existing style)
$ws.Range("B46").Style = $ws.Range("B45").Style

# Update selection to reflect new state (matching diff: activeCell/sqref now B50)
$ws.Range("B50").Select()
